{"js": "// The author's commit fixes a typo: \"PSS-IFF\" (double F) should read\n// \"PSS-IF\" (single F). There is exactly one occurrence of \"PSS-IFF\" in\n// the whole document, inside the \"Vorgehen\" section's final paragraph\n// (\"Abschlie\u00dfend werden Transformatoren ... mittels PSS-IFF modellierten\n// PSS ...\"). Everything else in the diff (style-id renames, proofErr\n// clean-up, divId removal, sdt reordering) is Word's own re-normalization\n// noise from opening/resaving the file and carries no semantic content\n// change, so we only need to perform this single text replacement.\n\nconst body = context.document.body;\nconst results = body.search(\"PSS-IFF\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"PSS-IF\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The author's commit fixes a typo: \"PSS-IFF\" (double F) should read\n# \"PSS-IF\" (single F). There is exactly one occurrence of \"PSS-IFF\" in\n# the whole document, inside the \"Vorgehen\" section's final paragraph\n# (\"Abschlie\u00dfend werden Transformatoren ... mittels PSS-IFF modellierten\n# PSS ...\"). Everything else in the diff (style-id renames, proofErr\n# clean-up, divId removal, sdt reordering) is Word's own re-normalization\n# noise from opening/resaving the file and carries no semantic content\n# change, so we only need to perform this single text replacement.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"PSS-IFF\"\n$find.Replacement.Text = \"PSS-IF\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
